$d = $word.ActiveDocument

# --- Step 1: insert all the new paragraphs (27/02/2018 diary entry) between
# the "function." paragraph and the bookmark paragraph, by inserting at a
# position strictly inside the "function." paragraph's text (just before its
# trailing pilcrow) so InsertXML appends new paragraphs after it instead of
# overwriting it or the following (bookmark) paragraph.
$pFunction = $d.Paragraphs.Item(22)
$rFunction = $pFunction.Range
$insertPoint = $d.Range($rFunction.End - 1, $rFunction.End - 1)

$mainXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p><w:r><w:t>27/02/2018</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Today </w:t></w:r><w:r><w:t>I have half sorted it</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> it now builds near the choke not exactly on the choke, but at least its progress</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">I think the issue was how the system uses the build sites i.e. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>StartingLocation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = 1, Natural = 2, Extension = 3 so what </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> did was add a fourth being choke, but I never used this I tried to subvert it using your select choke code, which the system wouldn't accept, so I looked to see how the initial build site was set and it was set using a simple Self().</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>GetStartLocation</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">(), so instead of going around the build site selection I added the following code to set the add and set the choke point to the built site </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>enum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        [</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ExecutableAction</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>("</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SelectChokeBuild</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>")]</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        public bool </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>SelectChokeBuild</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        {</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">           if (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Interface(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>).</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>buildingChoke</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> is </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TilePosition</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">           {</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                Interface(</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>).</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>baseLocations</w:t></w:r><w:proofErr w:type="spellEnd"/><w:proofErr w:type="gramEnd"/><w:r><w:t>[(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>BuildSite.Choke</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>] = Interface().</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>buildingChoke</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">                return </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>SwitchBuildToBase</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>((</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>BuildSite.Choke</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>);</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">           }</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">            return false;</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">            </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">        }</w:t></w:r></w:p><w:p><w:r><w:t>This is the code to get the choke point position</w:t></w:r></w:p><w:p><w:r><w:t>though I'm not entirely confident on how it works, I know that it does get the position of the chokepoint and converts it into a tile position</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">           </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>Interface(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>).</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>buildingChoke</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> = (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>targetChoke.getDistance</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">(new </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TilePosition</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>chokepoint.getSides</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">().second)) &lt; </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>targetChoke.getDistance</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">(new </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TilePosition</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>chokepoint.getSides</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">().first))) ? new </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TilePosition</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>chokepoint.getSides</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">().second) : new </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>TilePosition</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>chokepoint.getSides</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>().first);</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">after testing, </w:t></w:r><w:r><w:t>it now build on the correct side of the choke</w:t></w:r><w:r><w:t xml:space="preserve"> most of the time, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>i</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> basically swapped the </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>getsides</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:t xml:space="preserve"> as can be seen in the snippet</w:t></w:r><w:r><w:t>.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">In conclusion to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>todays</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> work, as</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>I</w:t></w:r><w:r><w:t xml:space="preserve"> said </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>originally</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>I</w:t></w:r><w:r><w:t xml:space="preserve"> was just passing that value straight into the get possible build location which it didn't like</w:t></w:r></w:p>
'@

$insertPoint.InsertXML($mainXml) | Out-Null

# --- Step 2: replace the (now shifted) trailing bookmark-only paragraph's
# content with the final "So my issue..." paragraph text, recreating the
# _GoBack bookmark so it still wraps an (empty) range at the very end.
$pBookmark = $d.Paragraphs.Last
$rBookmark = $pBookmark.Range
$rBookmark.Collapse(1)

$lastXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">So my issue again was just not understanding the </w:t></w:r><w:r><w:t>API.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$rBookmark.InsertXML($lastXml) | Out-Null

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
